{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target: the LinkedListTurbo bullet point. Before:\n//   \"indexedContains on datastructures.LinkedListTurbo is O(n2) because for\n//    contains it must walk through the entire list on worst case which is n\n//    and then the get method at worst case must also walk through the entire\n//    list which is another n making it O(n2)\"\n// (both \"2\"s are superscript runs)\n// After:\n//   \"indexedContains on datastructures.LinkedListTurbo is O(n) because for\n//    contains it must walk through the entire list on worst case which is n\n//    and then the get method each time will only walk 1 position and it will\n//    save its last position making it O(1) making the method overall O(n)\"\n// (no superscripts remain)\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Locate the paragraph that starts the \"indexedContains on\n// datastructures.LinkedListTurbo\" bullet (unique lead-in text).\nconst marker = \"indexedContains on datastructures.LinkedListTurbo\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(marker) === 0) {\n    target = para;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find target paragraph for LinkedListTurbo analysis.\");\n}\n\n// Search inside that paragraph for the lead-in run so we can grab the range\n// that spans everything after it (i.e. the two sentences plus both\n// superscript \"2\" runs that need to go away).\nconst hits = target.search(marker, { matchCase: true });\nhits.load(\"items\");\nawait context.sync();\n\nconst hit = hits.items[0];\nconst tail = hit.getRange(\"After\").expandTo(target.getRange(\"End\"));\n\n// Replace the whole tail (which currently contains the superscripted\n// \"O(n2)...O(n2)\" text split across several runs) with the new plain-text\n// explanation, written as two runs to mirror the authored edit.\ntail.insertText(\n  \" is O(n) because for contains it must walk through the entire list on worst case which is n and then the get method \",\n  \"Replace\"\n);\n\n// Re-resolve the paragraph/range after the mutation and append the closing\n// sentence at the (new) end of the paragraph.\nparagraphs.load(\"items\");\nawait context.sync();\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  para.load(\"text\");\n}\nawait context.sync();\n\nlet target2 = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  if (para.text.indexOf(marker) === 0) {\n    target2 = para;\n    break;\n  }\n}\n\nconst endRange = target2.getRange(\"End\");\nendRange.insertText(\n  \"each time will only walk 1 position and it will save its last position making it O(1) making the method overall O(n)\",\n  \"Start\"\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is open as $d.\n#\n# Target: the LinkedListTurbo bullet point. Before:\n#   \"indexedContains on datastructures.LinkedListTurbo is O(n2) because for\n#    contains it must walk through the entire list on worst case which is n\n#    and then the get method at worst case must also walk through the entire\n#    list which is another n making it O(n2)\"\n# (both \"2\"s are superscript runs)\n# After:\n#   \"indexedContains on datastructures.LinkedListTurbo is O(n) because for\n#    contains it must walk through the entire list on worst case which is n\n#    and then the get method each time will only walk 1 position and it will\n#    save its last position making it O(1) making the method overall O(n)\"\n# (no superscripts remain)\n\n$d = $word.ActiveDocument\n\n$marker = \"indexedContains on datastructures.LinkedListTurbo\"\n\n# Locate the paragraph that starts with the unique lead-in text (there is\n# an earlier heading that also contains \"LinkedListTurbo\", so match on the\n# full bullet lead-in rather than a bare substring search of the document).\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith($marker)) {\n        $target = $p\n        break\n    }\n}\n\nif ($target -eq $null) {\n    Write-Output \"Could not find target paragraph for LinkedListTurbo analysis.\"\n} else {\n    $pStart = $target.Range.Start\n    $pEnd = $target.Range.End\n\n    # Find the end of the \"indexedContains on datastructures.LinkedListTurbo\"\n    # lead-in within this paragraph only.\n    $findRng = $d.Range($pStart, $pEnd)\n    $found = $findRng.Find.Execute($marker)\n\n    # Range spanning everything after the lead-in through the end of the\n    # paragraph -- this currently holds \"O(n2)...O(n2)\" split across several\n    # runs (two of them superscripted \"2\"s). Replace it wholesale with the\n    # new plain-text explanation.\n    $tail = $d.Range($findRng.End, $pEnd)\n    $tail.Text = \" is O(n) because for contains it must walk through the entire list on worst case which is n and then the get method each time will only walk 1 position and it will save its last position making it O(1) making the method overall O(n)\"\n}\n"}
